# Rename the Treatment labels in column K ("Treatment") to their new,
# more descriptive names (header row 1 stays as "Treatment").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "FUT"    = "ERCP 8.5"
    "AMB"    = "Ambient"
    "AMB+HW" = "Ambient HW"
    "FUT+HW" = "ERCP 8.5 HW"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 2) {
    $lastRow = 193
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $current = $cell.Text
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
